# Append the 2025-12-26 Kaspa buy entry as a new row at the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A stores the date as literal text (matching the rest of the sheet's
# "MM/DD/YYYY" string entries), so force Text formatting before assigning the
# value to stop Excel from auto-converting it into a date serial number.
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "12/26/2025"
# Drop back to the default (unstyled) cell style now that the value is safely
# stored as text, so the new cell doesn't carry a stray explicit format.
$ws.Range("A20").Style = "Normal"

$ws.Range("B20").Value = 1120.191999999999
$ws.Range("C20").Value = 0.04418885333942756
$ws.Range("D20").Value = 50
